$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing all existing data down by one row
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "xcoord"
$ws.Range("C1").Value = "ycoord"

# Update the selected cell to A2
$ws.Range("A2").Select()
